$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.881.45'
$ws.Range('E2').Value = '  +2.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.884.57'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.12'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4676'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3943'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07929'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9836'
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.45'
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.933.65'
$ws.Range('E12').Value = '  +8.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.754'
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.021'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06986'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.94'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001012'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.01'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.862.40'
$ws.Range('E21').Value = '  +2.76%  '
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.12'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.121'
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.051.84'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.42'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.44'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.812'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.009'
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.11'
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09410'
$ws.Range('E31').Value = '  +1.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9435'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.323'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.363'
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.347'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05936'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02130'
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.164'
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('E39').Value = '  +3.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5739'
$ws.Range('E40').Value = '  +2.52%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.04'
$ws.Range('E41').Value = '  +1.20%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1802'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07327'
$ws.Range('E43').Value = '  +4.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.83'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.179'
$ws.Range('E45').Value = '  -3.71%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5361'
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.850'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.102'
$ws.Range('E48').Value = '  -5.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '114.26'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.381'
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('E51').Value = '  +0.60%  '
